$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" column header in H1 - copy the formatting used by the other
# header cells (bold, centered, bordered) then set its text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Corrected Prediction/Error values from refitting NCDEs to individual
# patients (rows 7,8,10,11 in block 1; block 2 values unchanged)
$ws.Range("D7").Value = 0.6006829552088109
$ws.Range("E7").Value = 0.3993170447911891
$ws.Range("D8").Value = 0.8026935976643621
$ws.Range("E8").Value = 0.1973064023356379
$ws.Range("D10").Value = 0.8503492011709716
$ws.Range("E10").Value = 0.1496507988290284
$ws.Range("D11").Value = 0.8310994405101909
$ws.Range("E11").Value = 0.1689005594898091

# Populate the new Label column: 0 for Control rows, 1 for MDD rows
$labels = @{
    2  = 0; 3  = 0; 4  = 0; 5  = 0; 6  = 0
    7  = 1; 8  = 1; 9  = 1; 10 = 1; 11 = 1
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0
    17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
